# ---------------------------------------------------------------------------
# Target change (from the supplied OOXML diff):
#
#   word/styles.xml -> <w:docDefaults> is trimmed down: every attribute that
#   merely restated the OOXML schema default (w:b=0, w:i=0, w:smallCaps=0,
#   w:strike=0, w:color=000000, w:u=none, w:shd=clear/auto, w:vertAlign=
#   baseline, w:keepNext=0, w:keepLines=0, w:widowControl=1, w:pBdr=nil,
#   w:spacing/@before=0,@after=0, w:ind=0, w:contextualSpacing=0, w:jc=left)
#   is removed, leaving only <w:rFonts>, <w:sz>, <w:szCs>, <w:lang> in
#   rPrDefault and a bare <w:spacing w:line="276" w:lineRule="auto"/> in
#   pPrDefault. The commit message ("download tc, tcn, and tl files from
#   GD") says this came from re-exporting the file from Google Drive, i.e.
#   it is a by-product of a whole-file regeneration pipeline, not an
#   in-app formatting edit - the <w:docDefaults> block is never touched by
#   the document's visible content (none of the removed attributes differ
#   from their implied defaults, and the one value that is kept,
#   line=276/auto, is left exactly as-is), so nothing about how the
#   document actually renders changes.
#
# Why this script does not attempt to rewrite <w:docDefaults>:
#   <w:docDefaults> is template-level markup that sits outside every
#   w:style element. Word's object model (both in real Word and in this
#   COM-interop surface) does not expose it: there is no
#   Document.DocDefaults / Styles.DocDefaults / Document.Defaults, etc.
#   - Styles("Normal").Font / .ParagraphFormat write directly into the
#     <w:style w:styleId="Normal"> element (which the diff shows must stay
#     exactly as it is: just <w:name w:val="normal"/>, no rPr/pPr), not
#     into docDefaults, so using it would add formatting the diff never
#     adds and would rewrite (re-serialize) every other style definition
#     along with it.
#   - Range/Content-level Font / ParagraphFormat assignment writes direct
#     run/paragraph formatting into word/document.xml on every paragraph,
#     which the diff shows is untouched.
#   - Find/Replace only ever touches document text content, never part-
#     level XML such as styles.xml.
#   - Round-tripping the document (Save / SaveAs2 / reopen) does not
#     normalize or rewrite docDefaults either.
#
# So there is no reachable Word-OM call that lands on <w:docDefaults>
# without corrupting a part of the package the diff leaves alone. Since
# the edit has no visible effect on the document (every dropped attribute
# already equalled its implied default, and the single surviving value is
# unchanged), the safe, faithful action is to touch the document only
# through read-only navigation - i.e. make no spurious mutation - rather
# than force a change through an API that does not correspond to this
# part of the file.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# Confirm the paragraph style that inherits from docDefaults is present;
# this is purely a read (no property is assigned) so it leaves the
# package byte-for-byte untouched, matching the fact that docDefaults is
# not reachable through this object model and nothing else in the
# document changes.
$normal = $d.Styles.Item("Normal")
$null = $normal.NameLocal
